# Workbook tracks daily % positivity (or similar) rates per US state.
# This edit:
#   1. Backfills American Samoa (column E) figures for 01-07 Jun 2020 (rows 123-129).
#   2. Completes state figures (columns B:BE) for 03-09 Sep 2020 (rows 217-223),
#      which previously only had the date label in column A.
#   3. Appends new date rows for 10-14 Sep 2020 (rows 224-228), date label only.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Backfill American Samoa (column E) for rows 123-129 ---
$ws.Range("E123").Value = 0.08824855699855701
$ws.Range("E124").Value = 0.066666666666667
$ws.Range("E125").Value = 0.07142857142857099
$ws.Range("E126").Value = 0.055555555555556
$ws.Range("E127").Value = 0.04
$ws.Range("E128").Value = 0.04
$ws.Range("E129").Value = 0.03448275862069

# --- 2 & 3. Complete rows 217-223 and append new date rows 224-228 ---
$ws.Range("B217").Value = 0.054051777181167
$ws.Range("C217").Value = 0.07438515689485201
$ws.Range("D217").Value = 0.0738604828078
$ws.Range("F217").Value = 0.0456171526476
$ws.Range("G217").Value = 0.058230047036429
$ws.Range("H217").Value = 0.06011620828467
$ws.Range("I217").Value = 0.054684060933052
$ws.Range("J217").Value = 0.044544280676471
$ws.Range("K217").Value = 0.049765153601888
$ws.Range("L217").Value = 0.059431951327948
$ws.Range("M217").Value = 0.06402341157105
$ws.Range("N217").Value = 0.05894148898075
$ws.Range("O217").Value = 0.033970797790813
$ws.Range("P217").Value = 0.090041035206897
$ws.Range("Q217").Value = 0.06672973413746799
$ws.Range("R217").Value = 0.056850158392668
$ws.Range("S217").Value = 0.068565728695483
$ws.Range("T217").Value = 0.072839862095621
$ws.Range("U217").Value = 0.064301138862924
$ws.Range("V217").Value = 0.069916863010095
$ws.Range("W217").Value = 0.054578443518419
$ws.Range("X217").Value = 0.050309643541666
$ws.Range("Y217").Value = 0.06187218647042
$ws.Range("Z217").Value = 0.055971056641306
$ws.Range("AA217").Value = 0.051994237079487
$ws.Range("AB217").Value = 0.074978637787797
$ws.Range("AC217").Value = 0.074528439225998
$ws.Range("AD217").Value = 0.073762473364187
$ws.Range("AE217").Value = 0.070223000239988
$ws.Range("AF217").Value = 0.058087651072401
$ws.Range("AG217").Value = 0.074703563831317
$ws.Range("AH217").Value = 0.087970894799088
$ws.Range("AI217").Value = 0.063224830834889
$ws.Range("AJ217").Value = 0.054005238369497
$ws.Range("AK217").Value = 0.056392354330879
$ws.Range("AL217").Value = 0.051842289603355
$ws.Range("AM217").Value = 0.054801511252347
$ws.Range("AN217").Value = 0.062533337142858
$ws.Range("AO217").Value = 0.06770515358792201
$ws.Range("AP217").Value = 0.047275609226473
$ws.Range("AQ217").Value = 0.058514156851789
$ws.Range("AR217").Value = 0.078790580284043
$ws.Range("AS217").Value = 0.054457336686406
$ws.Range("AT217").Value = 0.0634249737448
$ws.Range("AU217").Value = 0.07688945607342899
$ws.Range("AV217").Value = 0.06628315205109001
$ws.Range("AW217").Value = 0.057192161236652
$ws.Range("AX217").Value = 0.066037670806815
$ws.Range("AY217").Value = 0.053260287807968
$ws.Range("AZ217").Value = 0.08627631944175899
$ws.Range("BA217").Value = 0.059688520520676
$ws.Range("BB217").Value = 0.049873568594837
$ws.Range("BC217").Value = 0.066927199265464
$ws.Range("BD217").Value = 0.057483164472022
$ws.Range("BE217").Value = 0.06713071233567

$ws.Range("B218").Value = 0.042710485544877
$ws.Range("C218").Value = 0.061611162361701
$ws.Range("D218").Value = 0.059548700068907
$ws.Range("F218").Value = 0.036826053578075
$ws.Range("G218").Value = 0.052085677482251
$ws.Range("H218").Value = 0.056488130755608
$ws.Range("I218").Value = 0.059240254568204
$ws.Range("J218").Value = 0.057674723860616
$ws.Range("K218").Value = 0.047889740777501
$ws.Range("L218").Value = 0.056356195338835
$ws.Range("M218").Value = 0.048467729379673
$ws.Range("N218").Value = 0.054343223472767
$ws.Range("O218").Value = 0.027288756341795
$ws.Range("P218").Value = 0.063455833879765
$ws.Range("Q218").Value = 0.04476842660233
$ws.Range("R218").Value = 0.041047362260838
$ws.Range("S218").Value = 0.048114529755593
$ws.Range("T218").Value = 0.051572547524188
$ws.Range("U218").Value = 0.042167527062057
$ws.Range("V218").Value = 0.050090894039522
$ws.Range("W218").Value = 0.045485240086773
$ws.Range("X218").Value = 0.039700891826407
$ws.Range("Y218").Value = 0.047306242420121
$ws.Range("Z218").Value = 0.043224188540384
$ws.Range("AA218").Value = 0.043639501220807
$ws.Range("AB218").Value = 0.06321169344661901
$ws.Range("AC218").Value = 0.07555615362937899
$ws.Range("AD218").Value = 0.062149964924606
$ws.Range("AE218").Value = 0.061050345654124
$ws.Range("AF218").Value = 0.048997235883406
$ws.Range("AG218").Value = 0.06469543942195
$ws.Range("AH218").Value = 0.07388735108484901
$ws.Range("AI218").Value = 0.046330292073797
$ws.Range("AJ218").Value = 0.042553574061827
$ws.Range("AK218").Value = 0.047061575714342
$ws.Range("AL218").Value = 0.045570799809682
$ws.Range("AM218").Value = 0.048030477833557
$ws.Range("AN218").Value = 0.055898819229505
$ws.Range("AO218").Value = 0.056513541005481
$ws.Range("AP218").Value = 0.043897939364471
$ws.Range("AQ218").Value = 0.04640487720092
$ws.Range("AR218").Value = 0.074757237943558
$ws.Range("AS218").Value = 0.04530581810441
$ws.Range("AT218").Value = 0.04772673858807
$ws.Range("AU218").Value = 0.048935588185846
$ws.Range("AV218").Value = 0.051632382827469
$ws.Range("AW218").Value = 0.046244025070133
$ws.Range("AX218").Value = 0.051739318675116
$ws.Range("AY218").Value = 0.041465770376391
$ws.Range("AZ218").Value = 0.07650479673057101
$ws.Range("BA218").Value = 0.045924272586575
$ws.Range("BB218").Value = 0.041118400577428
$ws.Range("BC218").Value = 0.05498489602733
$ws.Range("BD218").Value = 0.046184119328981
$ws.Range("BE218").Value = 0.05241838759616

$ws.Range("B219").Value = 0.034381609421372
$ws.Range("C219").Value = 0.023613756377148
$ws.Range("D219").Value = 0.024007183724442
$ws.Range("F219").Value = 0.0272359926018
$ws.Range("G219").Value = 0.039932698939804
$ws.Range("H219").Value = 0.035669264853925
$ws.Range("I219").Value = 0.028364548033779
$ws.Range("J219").Value = 0.038813133925931
$ws.Range("K219").Value = 0.028293585281671
$ws.Range("L219").Value = 0.025847069505638
$ws.Range("M219").Value = 0.025633850225
$ws.Range("N219").Value = 0.059301279959332
$ws.Range("O219").Value = 0.025968482676658
$ws.Range("P219").Value = 0.026009347054522
$ws.Range("Q219").Value = 0.037314099130703
$ws.Range("R219").Value = 0.031769392314
$ws.Range("S219").Value = 0.026042966285571
$ws.Range("T219").Value = 0.024554261130003
$ws.Range("U219").Value = 0.02422300700049
$ws.Range("V219").Value = 0.029029932389594
$ws.Range("W219").Value = 0.033384501056168
$ws.Range("X219").Value = 0.027601208443583
$ws.Range("Y219").Value = 0.031435145021726
$ws.Range("Z219").Value = 0.027710414365475
$ws.Range("AA219").Value = 0.024419942959275
$ws.Range("AB219").Value = 0.023613539668638
$ws.Range("AC219").Value = 0.032258064516129
$ws.Range("AD219").Value = 0.022778905347017
$ws.Range("AE219").Value = 0.034243504832655
$ws.Range("AF219").Value = 0.025004043573436
$ws.Range("AG219").Value = 0.029315168864589
$ws.Range("AH219").Value = 0.024348215085731
$ws.Range("AI219").Value = 0.025585579596129
$ws.Range("AJ219").Value = 0.025878110140271
$ws.Range("AK219").Value = 0.034446816458073
$ws.Range("AL219").Value = 0.028530542622511
$ws.Range("AM219").Value = 0.032972284805662
$ws.Range("AN219").Value = 0.028893497987105
$ws.Range("AO219").Value = 0.025908583583811
$ws.Range("AP219").Value = 0.034229765970104
$ws.Range("AQ219").Value = 0.031909478738068
$ws.Range("AR219").Value = 0.071346429010763
$ws.Range("AS219").Value = 0.035831902468437
$ws.Range("AT219").Value = 0.026567885659019
$ws.Range("AU219").Value = 0.031277151712665
$ws.Range("AV219").Value = 0.025971984141687
$ws.Range("AW219").Value = 0.026724383787796
$ws.Range("AX219").Value = 0.036893379328037
$ws.Range("AY219").Value = 0.028943747983393
$ws.Range("AZ219").Value = 0.069499502555656
$ws.Range("BA219").Value = 0.038977036534405
$ws.Range("BB219").Value = 0.032738462403396
$ws.Range("BC219").Value = 0.028439819222274
$ws.Range("BD219").Value = 0.031070888494907
$ws.Range("BE219").Value = 0.03534211458212

$ws.Range("B220").Value = 0.043784947245174
$ws.Range("C220").Value = 0.025060070471132
$ws.Range("D220").Value = 0.021543243658368
$ws.Range("F220").Value = 0.02811157325704
$ws.Range("G220").Value = 0.028890323729104
$ws.Range("H220").Value = 0.029728421758617
$ws.Range("I220").Value = 0.026465502531631
$ws.Range("J220").Value = 0.040047894324592
$ws.Range("K220").Value = 0.025367020618753
$ws.Range("L220").Value = 0.023555198390088
$ws.Range("M220").Value = 0.021362822847995
$ws.Range("N220").Value = 0.055994557088414
$ws.Range("O220").Value = 0.024565848184243
$ws.Range("P220").Value = 0.02318692058494
$ws.Range("Q220").Value = 0.031041903302498
$ws.Range("R220").Value = 0.026303239043875
$ws.Range("S220").Value = 0.023250156828383
$ws.Range("T220").Value = 0.023704754794478
$ws.Range("U220").Value = 0.022783877853218
$ws.Range("V220").Value = 0.026779133411092
$ws.Range("W220").Value = 0.030515264762792
$ws.Range("X220").Value = 0.025433332794019
$ws.Range("Y220").Value = 0.029847800865596
$ws.Range("Z220").Value = 0.025354793579908
$ws.Range("AA220").Value = 0.02331731861083
$ws.Range("AB220").Value = 0.022075355179746
$ws.Range("AC220").Value = 0.045454545454545
$ws.Range("AD220").Value = 0.020402645857307
$ws.Range("AE220").Value = 0.033835417805526
$ws.Range("AF220").Value = 0.023526530900028
$ws.Range("AG220").Value = 0.027372239421361
$ws.Range("AH220").Value = 0.023601990737471
$ws.Range("AI220").Value = 0.025485590047435
$ws.Range("AJ220").Value = 0.026805562617296
$ws.Range("AK220").Value = 0.037645384679056
$ws.Range("AL220").Value = 0.02769980135361
$ws.Range("AM220").Value = 0.035696861679525
$ws.Range("AN220").Value = 0.028133588969447
$ws.Range("AO220").Value = 0.023061134286465
$ws.Range("AP220").Value = 0.030611739458278
$ws.Range("AQ220").Value = 0.028249586634902
$ws.Range("AR220").Value = 0.066481547104718
$ws.Range("AS220").Value = 0.031704695510383
$ws.Range("AT220").Value = 0.020777282822247
$ws.Range("AU220").Value = 0.025097244647127
$ws.Range("AV220").Value = 0.019773605151081
$ws.Range("AW220").Value = 0.019605830460213
$ws.Range("AX220").Value = 0.026984931923115
$ws.Range("AY220").Value = 0.022962814583392
$ws.Range("AZ220").Value = 0.06186885577657
$ws.Range("BA220").Value = 0.032640877488157
$ws.Range("BB220").Value = 0.026924183217225
$ws.Range("BC220").Value = 0.025143614284906
$ws.Range("BD220").Value = 0.026755561924353
$ws.Range("BE220").Value = 0.02931026919049

$ws.Range("B221").Value = 0.043980847512857
$ws.Range("C221").Value = 0.022954720393222
$ws.Range("D221").Value = 0.019679404317117
$ws.Range("F221").Value = 0.026625982911938
$ws.Range("G221").Value = 0.033221269133641
$ws.Range("H221").Value = 0.030452139519365
$ws.Range("I221").Value = 0.029242945369842
$ws.Range("J221").Value = 0.048047812661014
$ws.Range("K221").Value = 0.027837150759563
$ws.Range("L221").Value = 0.029073596980399
$ws.Range("M221").Value = 0.025377834067554
$ws.Range("N221").Value = 0.053273545902358
$ws.Range("O221").Value = 0.026096423860116
$ws.Range("P221").Value = 0.021101320371908
$ws.Range("Q221").Value = 0.032999612131734
$ws.Range("R221").Value = 0.02718167379395
$ws.Range("S221").Value = 0.021815001171181
$ws.Range("T221").Value = 0.021368080590068
$ws.Range("U221").Value = 0.021797304122306
$ws.Range("V221").Value = 0.027716425392725
$ws.Range("W221").Value = 0.032045401178898
$ws.Range("X221").Value = 0.026708173818026
$ws.Range("Y221").Value = 0.029340137932684
$ws.Range("Z221").Value = 0.026121805912437
$ws.Range("AA221").Value = 0.024201528248643
$ws.Range("AB221").Value = 0.02389600769207
$ws.Range("AC221").Value = 0.10714285714286
$ws.Range("AD221").Value = 0.021826406642218
$ws.Range("AE221").Value = 0.035496701708366
$ws.Range("AF221").Value = 0.027038752625036
$ws.Range("AG221").Value = 0.029960792259418
$ws.Range("AH221").Value = 0.024688915202124
$ws.Range("AI221").Value = 0.028797541604559
$ws.Range("AJ221").Value = 0.030898445692022
$ws.Range("AK221").Value = 0.040348641208058
$ws.Range("AL221").Value = 0.030476117476802
$ws.Range("AM221").Value = 0.042203957674751
$ws.Range("AN221").Value = 0.027724594173255
$ws.Range("AO221").Value = 0.0237141537977
$ws.Range("AP221").Value = 0.029762343978024
$ws.Range("AQ221").Value = 0.029078951039638
$ws.Range("AR221").Value = 0.06737068830863201
$ws.Range("AS221").Value = 0.032774003118454
$ws.Range("AT221").Value = 0.024185487521417
$ws.Range("AU221").Value = 0.026623606837623
$ws.Range("AV221").Value = 0.022812960829811
$ws.Range("AW221").Value = 0.023671292738487
$ws.Range("AX221").Value = 0.029766220339437
$ws.Range("AY221").Value = 0.024988737718337
$ws.Range("AZ221").Value = 0.069773444283035
$ws.Range("BA221").Value = 0.035199854496327
$ws.Range("BB221").Value = 0.027881825286638
$ws.Range("BC221").Value = 0.024359995389024
$ws.Range("BD221").Value = 0.027372931945649
$ws.Range("BE221").Value = 0.028578037205865

$ws.Range("A222").Value = "08 09 2020"
$ws.Range("B222").Value = 0.055099777617582
$ws.Range("C222").Value = 0.07357331436457699
$ws.Range("D222").Value = 0.067162863911504
$ws.Range("F222").Value = 0.042227893729904
$ws.Range("G222").Value = 0.055604700738256
$ws.Range("H222").Value = 0.062664452205533
$ws.Range("I222").Value = 0.060856992205913
$ws.Range("J222").Value = 0.048465448760751
$ws.Range("K222").Value = 0.052948754986868
$ws.Range("L222").Value = 0.060485247576582
$ws.Range("M222").Value = 0.06491697745646401
$ws.Range("N222").Value = 0.064811260357926
$ws.Range("O222").Value = 0.03392480352676
$ws.Range("P222").Value = 0.09241499059403201
$ws.Range("Q222").Value = 0.066863083530531
$ws.Range("R222").Value = 0.05913373606542
$ws.Range("S222").Value = 0.06734183708540301
$ws.Range("T222").Value = 0.073259243456368
$ws.Range("U222").Value = 0.063371356321845
$ws.Range("V222").Value = 0.072209276921379
$ws.Range("W222").Value = 0.057370186527555
$ws.Range("X222").Value = 0.052632919206198
$ws.Range("Y222").Value = 0.065147063612305
$ws.Range("Z222").Value = 0.06291834134391799
$ws.Range("AA222").Value = 0.063714079550789
$ws.Range("AB222").Value = 0.07742894282795899
$ws.Range("AC222").Value = 0.10714285714286
$ws.Range("AD222").Value = 0.07786913965181801
$ws.Range("AE222").Value = 0.07610031323145799
$ws.Range("AF222").Value = 0.058468199294164
$ws.Range("AG222").Value = 0.076198249162045
$ws.Range("AH222").Value = 0.088551008539547
$ws.Range("AI222").Value = 0.059899615561033
$ws.Range("AJ222").Value = 0.05196297807995
$ws.Range("AK222").Value = 0.054068569092454
$ws.Range("AL222").Value = 0.050744920457597
$ws.Range("AM222").Value = 0.056265134250658
$ws.Range("AN222").Value = 0.06464039040494
$ws.Range("AO222").Value = 0.06756415696542099
$ws.Range("AP222").Value = 0.046484356819168
$ws.Range("AQ222").Value = 0.058024066420828
$ws.Range("AR222").Value = 0.071298492035621
$ws.Range("AS222").Value = 0.05206948716924
$ws.Range("AT222").Value = 0.063885672041294
$ws.Range("AU222").Value = 0.08221328976584
$ws.Range("AV222").Value = 0.065989722872595
$ws.Range("AW222").Value = 0.061250006148402
$ws.Range("AX222").Value = 0.065497474822725
$ws.Range("AY222").Value = 0.060743709960693
$ws.Range("AZ222").Value = 0.058522305405861
$ws.Range("BA222").Value = 0.06760100184730899
$ws.Range("BB222").Value = 0.054115320275422
$ws.Range("BC222").Value = 0.07419223988535199
$ws.Range("BD222").Value = 0.067104246970036
$ws.Range("BE222").Value = 0.069954580806439

$ws.Range("A223").Value = "09 09 2020"
$ws.Range("B223").Value = 0.06233703029946
$ws.Range("C223").Value = 0.08006736483764
$ws.Range("D223").Value = 0.078665537607293
$ws.Range("F223").Value = 0.05145993946714
$ws.Range("G223").Value = 0.05599263576615
$ws.Range("H223").Value = 0.056727521657554
$ws.Range("I223").Value = 0.051817614562139
$ws.Range("J223").Value = 0.042324555092455
$ws.Range("K223").Value = 0.046257925958561
$ws.Range("L223").Value = 0.051849344879533
$ws.Range("M223").Value = 0.057685094193455
$ws.Range("N223").Value = 0.055589774328391
$ws.Range("O223").Value = 0.027888675402088
$ws.Range("P223").Value = 0.07997875443718
$ws.Range("Q223").Value = 0.060472466891424
$ws.Range("R223").Value = 0.051360043981396
$ws.Range("S223").Value = 0.060969764649746
$ws.Range("T223").Value = 0.07081561967786
$ws.Range("U223").Value = 0.05642441603209
$ws.Range("V223").Value = 0.06540190937344501
$ws.Range("W223").Value = 0.046814722028883
$ws.Range("X223").Value = 0.044338872092449
$ws.Range("Y223").Value = 0.053933634625667
$ws.Range("Z223").Value = 0.054131186289087
$ws.Range("AA223").Value = 0.05470442484022
$ws.Range("AB223").Value = 0.071476856941772
$ws.Range("AC223").Value = 0.032258064516129
$ws.Range("AD223").Value = 0.070418458535945
$ws.Range("AE223").Value = 0.06827130011422899
$ws.Range("AF223").Value = 0.055671087878886
$ws.Range("AG223").Value = 0.0746521851787
$ws.Range("AH223").Value = 0.085599096141524
$ws.Range("AI223").Value = 0.059580211976321
$ws.Range("AJ223").Value = 0.052608317022114
$ws.Range("AK223").Value = 0.051989302263832
$ws.Range("AL223").Value = 0.047567007410545
$ws.Range("AM223").Value = 0.05681831632375
$ws.Range("AN223").Value = 0.061177102746296
$ws.Range("AO223").Value = 0.065900106923635
$ws.Range("AP223").Value = 0.041494572592151
$ws.Range("AQ223").Value = 0.054911192163063
$ws.Range("AR223").Value = 0.065356939548303
$ws.Range("AS223").Value = 0.051009241524045
$ws.Range("AT223").Value = 0.060710762661969
$ws.Range("AU223").Value = 0.07792611923352299
$ws.Range("AV223").Value = 0.062398764557326
$ws.Range("AW223").Value = 0.057686366700999
$ws.Range("AX223").Value = 0.062039264807687
$ws.Range("AY223").Value = 0.052831348530397
$ws.Range("AZ223").Value = 0.06256067939188301
$ws.Range("BA223").Value = 0.053923497487543
$ws.Range("BB223").Value = 0.045992159899892
$ws.Range("BC223").Value = 0.064230531404093
$ws.Range("BD223").Value = 0.057844411081379
$ws.Range("BE223").Value = 0.06479847427334801

$ws.Range("A224").Value = "10 09 2020"
$ws.Range("A225").Value = "11 09 2020"
$ws.Range("A226").Value = "12 09 2020"
$ws.Range("A227").Value = "13 09 2020"
$ws.Range("A228").Value = "14 09 2020"
